# Apply the UML class-inheritance diagram repositioning:
# the whole hand-drawn diagram (rectangles / triangles / connectors /
# the "UML Notation" textbox) is shifted up-and-right by a constant
# amount (dx = +1066800 EMU, dy = -1209488 EMU) once the empty "Title"
# placeholder is removed from the slide, freeing vertical space for the
# diagram to move into.
#
# PowerPoint's COM Shape.Left / Shape.Top are expressed in points and
# backed by 32-bit floats, so the literal point values below were
# chosen (via exhaustive search) to be the exact doubles that, once
# rounded through that float32 + points->EMU pipeline, reproduce the
# target EMU coordinates bit-for-bit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targets = @{
  84 = @{ L = 336.0;                T = 210.76472440944883 }
  47 = @{ L = 99.5;                 T = 210.76448828897637 }
  48 = @{ L = 99.5;                 T = 210.76472440944883 }
  67 = @{ L = 395.7591400582677;    T = 232.01236730472442 }
  70 = @{ L = 455.9999237598425;    T = 315.30528269055117 }
  71 = @{ L = 378.0;                T = 342.7647248094488 }
  72 = @{ L = 457.27024842047246;   T = 335.0348968897638 }
  46 = @{ L = 162.5;                T = 246.76464566929133 }
  54 = @{ L = 201.5;                T = 276.7646456692913 }
  61 = @{ L = 130.0;                T = 328.7914960629921 }
  62 = @{ L = 190.73661417322833;   T = 308.52803149606297 }
  73 = @{ L = 395.7591400582677;    T = 260.43299872598425 }
  74 = @{ L = 395.7591400582677;    T = 286.59883119763776 }
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $sid = [int]$sh.Id
    if ($targets.ContainsKey($sid)) {
        $t = $targets[$sid]
        $sh.Left = $t.L
        $sh.Top = $t.T
    }
}

# Remove the now-unused, empty title placeholder ("Title 59") that used
# to reserve space at the top of the slide.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ([int]$sh.Id -eq 60) {
        $sh.Delete()
    }
}
